$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bulunmuyor = "urun hafele.com.tr de bulunmuyor"

# Remove rows 5-11 (no longer part of the result set)
$ws.Range("A5:F11").EntireRow.Delete()

# Row 2: 900.70.391, all columns say "urun hafele.com.tr de bulunmuyor"
$ws.Range("A2").Value = "900.70.391"
$ws.Range("B2").Value = $bulunmuyor
$ws.Range("C2").Value = $bulunmuyor
$ws.Range("D2").Value = $bulunmuyor
$ws.Range("E2").Value = $bulunmuyor
$ws.Range("F2").Value = $bulunmuyor

# Row 3: 900.78.217, all columns say "urun hafele.com.tr de bulunmuyor"
$ws.Range("A3").Value = "900.78.217"
$ws.Range("B3").Value = $bulunmuyor
$ws.Range("C3").Value = $bulunmuyor
$ws.Range("D3").Value = $bulunmuyor
$ws.Range("E3").Value = $bulunmuyor
$ws.Range("F3").Value = $bulunmuyor

# Row 4: 900.78.417, with numeric stock and set urun pricing
$ws.Range("A4").Value = "900.78.417"
$ws.Range("B4").Value = 1238
$ws.Range("C4").Value = "set urun"
$ws.Range("D4").Value = "2.040,34 TL"
$ws.Range("E4").Value = "1.360,22 TL"
$ws.Range("F4").Value = "1.768,29 TL"
